$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AP column: replace the power-law slope formula with a simple ratio AL/AK ---
# Previously AP2 was blank, AP3 had its own (unshared) formula, AP4:AP18 was a shared formula.
# Now AP2:AP18 all use "=AL{row}/AK{row}", with AP2/AP3 written individually (no fill style)
# and AP4:AP18 remaining a shared formula group.
$ws.Range("AP3").Style = "Normal"
$ws.Range("AP4").Style = "Normal"
$ws.Range("AP2").Formula = "=AL2/AK2"
$ws.Range("AP3").Formula = "=AL3/AK3"
$ws.Range("AP4:AP18").Formula = "=AL4/AK4"

# --- AM column: switch the reference from AJ to AI ---
# Previously AM2:AM18 was one shared formula group referencing AJ.
# Now AM2/AM3 get individual formulas referencing AI, while AM4:AM18 becomes its own
# shared formula group (also referencing AI).
$ws.Range("AM4:AM18").Formula = "=(-2*AI4-3)/(-2*AI4+6)"
$ws.Range("AM3").Formula = "=(-2*AI3-3)/(-2*AI3+6)"
$ws.Range("AM2").Formula = "=(-2*AI2-3)/(-2*AI2+6)"

# --- Selection: whole column AM selected, active cell AM1 ---
$ws.Range("AM1:AM1048576").Select()

$wb.Save()
